$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2015503875968992
$ws.Range("C2").Value = 0.5542635658914729
$ws.Range("J2").Value = 0.01162790697674419
$ws.Range("P2").Value = 0.1589147286821705
$ws.Range("S2").Value = 0.07364341085271318
$ws.Range("B3").Value = 0.02068965517241379
$ws.Range("C3").Value = 0.01379310344827586
$ws.Range("J3").Value = 0.02758620689655172
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.1379310344827586
$ws.Range("J4").Value = 0.05357142857142857
$ws.Range("O4").Value = 0.01785714285714286
$ws.Range("P4").Value = 0.625
$ws.Range("S4").Value = 0.3035714285714285
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.04761904761904762
$ws.Range("D6").Value = 0.02645502645502645
$ws.Range("F6").Value = 0.07936507936507936
$ws.Range("J6").Value = 0.2380952380952381
$ws.Range("O6").Value = 0.02116402116402116
$ws.Range("Q6").Value = 0.2275132275132275
$ws.Range("S6").Value = 0.3121693121693122
$ws.Range("B7").Value = 0.09392265193370165
$ws.Range("D7").Value = 0.03867403314917127
$ws.Range("E7").Value = 0.005524861878453038
$ws.Range("F7").Value = 0.06077348066298342
$ws.Range("J7").Value = 0.08839779005524862
$ws.Range("O7").Value = 0.02209944751381215
$ws.Range("Q7").Value = 0.2044198895027624
$ws.Range("R7").Value = 0.09944751381215469
$ws.Range("S7").Value = 0.3867403314917127
$ws.Range("B8").Value = 0.103448275862069
$ws.Range("D8").Value = 0.01622718052738337
$ws.Range("E8").Value = 0.002028397565922921
$ws.Range("F8").Value = 0.06288032454361055
$ws.Range("J8").Value = 0.1054766734279919
$ws.Range("O8").Value = 0.01419878296146045
$ws.Range("Q8").Value = 0.2251521298174442
$ws.Range("R8").Value = 0.08316430020283976
$ws.Range("S8").Value = 0.3874239350912779
$ws.Range("B9").Value = 0.1100917431192661
$ws.Range("D9").Value = 0.01834862385321101
$ws.Range("F9").Value = 0.06422018348623854
$ws.Range("J9").Value = 0.06880733944954129
$ws.Range("O9").Value = 0.03669724770642202
$ws.Range("Q9").Value = 0.1834862385321101
$ws.Range("R9").Value = 0.07798165137614679
$ws.Range("S9").Value = 0.4403669724770642
$ws.Range("B10").Value = 0.0845771144278607
$ws.Range("D10").Value = 0.03150912106135987
$ws.Range("E10").Value = 0.001658374792703151
$ws.Range("F10").Value = 0.04975124378109453
$ws.Range("J10").Value = 0.1243781094527363
$ws.Range("O10").Value = 0.02404643449419569
$ws.Range("Q10").Value = 0.2222222222222222
$ws.Range("R10").Value = 0.08374792703150911
$ws.Range("S10").Value = 0.3781094527363184
$ws.Range("G11").Value = 0.1312056737588652
$ws.Range("J11").Value = 0.09574468085106383
$ws.Range("K11").Value = 0.1843971631205674
$ws.Range("L11").Value = 0.5567375886524822
$ws.Range("S11").Value = 0.03191489361702127
$ws.Range("G12").Value = 0.7636363636363637
$ws.Range("J12").Value = 0.1636363636363636
$ws.Range("K12").Value = 0.006060606060606061
$ws.Range("L12").Value = 0.0303030303030303
$ws.Range("S12").Value = 0.03636363636363636
$ws.Range("G13").Value = 0.6097560975609756
$ws.Range("J13").Value = 0.3414634146341464
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.01731601731601732
$ws.Range("H15").Value = 0.1948051948051948
$ws.Range("I15").Value = 0.04761904761904762
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.09090909090909091
$ws.Range("M15").Value = 0.004329004329004329
$ws.Range("N15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.09090909090909091
$ws.Range("S15").Value = 0.2164502164502164
$ws.Range("F16").Value = 0.005405405405405406
$ws.Range("H16").Value = 0.2108108108108108
$ws.Range("I16").Value = 0.07027027027027027
$ws.Range("J16").Value = 0.4324324324324325
$ws.Range("K16").Value = 0.1297297297297297
$ws.Range("M16").Value = 0.02162162162162162
$ws.Range("O16").Value = 0.05405405405405406
$ws.Range("S16").Value = 0.07567567567567568
$ws.Range("F17").Value = 0.01004016064257028
$ws.Range("H17").Value = 0.2309236947791165
$ws.Range("I17").Value = 0.09839357429718876
$ws.Range("J17").Value = 0.3755020080321285
$ws.Range("K17").Value = 0.0823293172690763
$ws.Range("M17").Value = 0.02610441767068273
$ws.Range("O17").Value = 0.0783132530120482
$ws.Range("S17").Value = 0.09839357429718876
$ws.Range("F18").Value = 0.01639344262295082
$ws.Range("H18").Value = 0.1530054644808743
$ws.Range("I18").Value = 0.1147540983606557
$ws.Range("J18").Value = 0.459016393442623
$ws.Range("K18").Value = 0.08196721311475409
$ws.Range("M18").Value = 0.02185792349726776
$ws.Range("O18").Value = 0.08196721311475409
$ws.Range("S18").Value = 0.07103825136612021
$ws.Range("F19").Value = 0.0151006711409396
$ws.Range("H19").Value = 0.2197986577181208
$ws.Range("I19").Value = 0.1023489932885906
$ws.Range("J19").Value = 0.3783557046979866
$ws.Range("K19").Value = 0.1006711409395973
$ws.Range("M19").Value = 0.01593959731543624
$ws.Range("N19").Value = 0.0008389261744966443
$ws.Range("O19").Value = 0.05453020134228188
$ws.Range("S19").Value = 0.1124161073825503

Write-Output "Applied 114 cell updates"
